# RTM.xlsx update: map new CYRS & HSI deliveries to the CRS rows, and
# refresh the HSI requirement versions (V1.0 -> V1.2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-shape the merged ranges in columns C:D (CYRS column).
#    Previously CYRS cells spanned several rows (C5:D6, C7:D10, C11:D12)
#    because most rows had no CYRS mapping yet. Now every CRS row gets
#    its own CYRS requirement, so each C:D pair becomes its own merge.
# ---------------------------------------------------------------------
$ws.Range("C5:D6").UnMerge()
$ws.Range("C7:D10").UnMerge()
$ws.Range("C11:D12").UnMerge()

$ws.Range("C5:D5").Merge()
$ws.Range("C6:D6").Merge()
$ws.Range("C7:D7").Merge()
$ws.Range("C8:D8").Merge()
$ws.Range("C9:D9").Merge()
$ws.Range("C10:D10").Merge()
$ws.Range("C11:D11").Merge()
$ws.Range("C12:D12").Merge()

# ---------------------------------------------------------------------
# 2) Update the CRS text in column B (re-ordering / re-wording some
#    entries, trimming trailing spaces, merging CRS_04 & CRS_06 into a
#    single row 8).
# ---------------------------------------------------------------------
$ws.Range("B4").Value = "[DIGELV _CRS_01_V1.0]"
$ws.Range("B7").Value = "[DIGELV _CRS_05_V1.0]"
$ws.Range("B8").Value = "[DIGELV _CRS_04_V1.0]                                                [DIGELV _CRS_06_V1.0]"
$ws.Range("B9").Value = "[DIGELV _CRS_09_V1.0]"
$ws.Range("B10").Value = "[DIGELV _CRS_08_V1.0]"
$ws.Range("B11").Value = "[DIGELV _CRS_07_V1.0]"
$ws.Range("B12").Value = "[DIGELV _CRS_07_V1.0]  "

# ---------------------------------------------------------------------
# 3) Fill in the new CYRS requirement mapping in column C (one CYRS
#    requirement per CRS row now, each flagged IMP(SW)/IMP(HW)).
# ---------------------------------------------------------------------
$ws.Range("C4").Value = "Req_DIGELV_CYRS_01_V1.0_IMP(SW)"
$ws.Range("C5").Value = "Req_DIGELV_CYRS_02_V1.0_IMP(SW)"
$ws.Range("C6").Value = "Req_DIGELV_CYRS_03_V1.0_IMP(SW)"
$ws.Range("C7").Value = "Req_DIGELV_CYRS_04_V1.0_IMP(SW)"
$ws.Range("C8").Value = "Req_DIGELV_CYRS_05_V1.0_IMP(HW)"
$ws.Range("C9").Value = "Req_DIGELV_CYRS_06_V1.0_IMP(HW)"
$ws.Range("C10").Value = "Req_DIGELV_CYRS_07_V1.0_IMP(SW)"
$ws.Range("C11").Value = "Req_DIGELV_CYRS_08_V1.0_IMP(HW)"
$ws.Range("C12").Value = "Req_DIGELV_CYRS_09_V1.0_IMP(HW)"
$ws.Range("C13").Value = "Req_DIGELV_CYRS_010_V1.0_IMP(SW)"

# ---------------------------------------------------------------------
# 4) Bump the HSI requirement versions from V1.0 to V1.2.
# ---------------------------------------------------------------------
$ws.Range("G6").Value = "Req _ DIGELV _HSI_01_V1.2"
$ws.Range("G7").Value = "Req _ DIGELV _ HSI _02_V1.2"
$ws.Range("G11").Value = "Req _ DIGELV _ HSI _03_V1.2"
$ws.Range("G13").Value = "Req _ DIGELV _ HSI _04_V1.2"

# ---------------------------------------------------------------------
# 5) Row heights: several rows need to grow to fit the new wrapped
#    CYRS/CRS text.
# ---------------------------------------------------------------------
$ws.Rows(5).RowHeight = 15.75
$ws.Rows(8).RowHeight = 63
$ws.Rows(9).RowHeight = 63
$ws.Rows(10).RowHeight = 47.25
$ws.Rows(12).RowHeight = 63

# ---------------------------------------------------------------------
# 6) Column widths: B shrinks (no longer needs bestFit-wide text), D
#    grows (now holds wrapped CYRS requirement text).
# ---------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 51.3
$ws.Columns("D").ColumnWidth = 25.8

# ---------------------------------------------------------------------
# 7) Scroll / selection: the author last left the sheet scrolled right
#    with cell I10 selected.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("I10").Select()
